# Update 北京-漫展信息.xlsx per gh-pages data refresh (commit 456a3b4)
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 1299
$ws.Range("G6").Value = 118
$ws.Range("F7").Value = 383
$ws.Range("F8").Value = 8351
$ws.Range("F10").Value = 10256
$ws.Range("F21").Value = 69
$ws.Range("F23").Value = 398
$ws.Range("F25").Value = 1754
$ws.Range("F27").Value = 516
$ws.Range("F30").Value = 52
$ws.Range("F31").Value = 569
$ws.Range("F33").Value = 1080
$ws.Range("F35").Value = 40
$ws.Range("F36").Value = 1407
$ws.Range("F40").Value = 14
$ws.Range("F44").Value = 74
$ws.Range("F48").Value = 60
$ws.Range("F49").Value = 62

# --- Sheet "演出" ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("G4").Value = 108
$ws.Range("F19").Value = 371

# --- Sheet "本地生活" ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 2781

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F9").Value = 1299
$ws.Range("G11").Value = 118
$ws.Range("G12").Value = 108
$ws.Range("F13").Value = 8351
$ws.Range("F15").Value = 10257
$ws.Range("F21").Value = 398
$ws.Range("F22").Value = 1754
$ws.Range("F25").Value = 52
$ws.Range("F27").Value = 569
$ws.Range("F30").Value = 40
$ws.Range("F33").Value = 1407
$ws.Range("F41").Value = 74
$ws.Range("F45").Value = 371
$ws.Range("F48").Value = 60
$ws.Range("F49").Value = 62
